$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = -11.89109999999999
$ws.Range("E4").Value = 12.8732
$ws.Range("E6").Value = 12.1092
$ws.Range("A9").Value = -20.07619999999998
$ws.Range("E10").Value = 11.76239999999999
$ws.Range("C11").Value = -13.9614
$ws.Range("E11").Value = 13.34629999999999
$ws.Range("A18").Value = -23.07320000000001
$ws.Range("A20").Value = -22.03260000000002
$ws.Range("D21").Value = -7.303400000000002
$ws.Range("E21").Value = 13.50520000000001
